$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.145.39"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "'1.815.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").Value = "'233.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.32%  "
$ws.Range("D6").Value = "'0.613"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").Value = "'41.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("E9").Value = "  +6.31%  "
$ws.Range("D10").Value = "'0.0685"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "'2.077.92"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "'1.814.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").Value = "'11.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("D15").Value = "'0.660"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "'4.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "'35.083.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'69.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'0.0₃0792"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'239.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "
$ws.Range("D21").Value = "'11.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D25").Value = "'173.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "'7.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  -2.41%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  +18.92%  "
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("D32").Value = "'3.331.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E35").Value = "  -7.46%  "
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").Value = "'92.88"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.41%  "
$ws.Range("D38").Value = "'0.682"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'1.308.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.49%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'14.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  -5.33%  "
$ws.Range("D46").Value = "'2.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("E47").Value = "  +4.59%  "
$ws.Range("E48").Value = "  -1.16%  "
$ws.Range("D49").Value = "'1.992.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "'0.0648"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.84%  "
